# Updated cryptos list with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# coin table, and fixes the ordering of the Stacks / FirstDigitalUSD rows.
#
# Note: several "Price" values look numeric (e.g. 603.14) but the sheet
# stores them as text (t="inlineStr"/shared string, not t="n"). Assigning a
# plain numeric-looking string to Range.Value lets Excel auto-convert it to
# a real number, so for those cells we prefix the value with a leading
# apostrophe to force Excel to keep (store) it as text, exactly like a user
# typing '603.14 into the cell would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Forces the cell to hold $text as a text value, even if it looks numeric.
    $ws.Range($cellRef).Value = "'" + $text
}

# Row 2 - Bitcoin
$ws.Range('D2').Value = '65.810.50'
$ws.Range('E2').Value = '  +0.42%  '

# Row 3 - Ethereum
$ws.Range('D3').Value = '2.678.84'
$ws.Range('E3').Value = '  +0.68%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  +0.05%  '

# Row 5 - BNB
Set-TextValue 'D5' '603.14'
$ws.Range('E5').Value = '  -0.39%  '

# Row 6 - Solana
Set-TextValue 'D6' '156.47'
$ws.Range('E6').Value = '  -1.01%  '

# Row 7 - USDC
$ws.Range('E7').Value = '  +0.05%  '

# Row 8 - XRP
Set-TextValue 'D8' '0.590'
$ws.Range('E8').Value = '  +0.05%  '

# Row 9 - Dogecoin
$ws.Range('E9').Value = '  -0.11%  '

# Row 10 - Toncoin
Set-TextValue 'D10' '5.95'
$ws.Range('E10').Value = '  +2.09%  '

# Row 11 - Cardano
$ws.Range('E11').Value = '  -2.86%  '

# Row 12 - TRON
$ws.Range('E12').Value = '  +0.31%  '

# Row 13 - Avalanche
Set-TextValue 'D13' '29.44'
$ws.Range('E13').Value = '  -0.75%  '

# Row 14 - ShibaInu
$ws.Range('E14').Value = '  +5.76%  '

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '3.162.03'

# Row 16 - WrappedBTC
$ws.Range('D16').Value = '65.601.72'
$ws.Range('E16').Value = '  +0.52%  '

# Row 17 - WrappedEther
$ws.Range('D17').Value = '2.684.39'
$ws.Range('E17').Value = '  +1.28%  '

# Row 18 - Chainlink
Set-TextValue 'D18' '12.61'
$ws.Range('E18').Value = '  -1.64%  '

# Row 19 - Polkadot
$ws.Range('E19').Value = '  -1.87%  '

# Row 20 - Uniswap
Set-TextValue 'D20' '7.58'
$ws.Range('E20').Value = '  +2.87%  '

# Row 21 - BitcoinCash
Set-TextValue 'D21' '352.07'
$ws.Range('E21').Value = '  -2.34%  '

# Row 22 - Dai
$ws.Range('E22').Value = '  -0.12%  '

# Row 23 - Litecoin
Set-TextValue 'D23' '70.43'
$ws.Range('E23').Value = '  +1.87%  '

# Row 24 - PEPE
Set-TextValue 'D24' '0.0000110'
$ws.Range('E24').Value = '  +6.24%  '

# Row 25 - InternetComputer(DFINITY)
Set-TextValue 'D25' '9.82'
$ws.Range('E25').Value = '  +2.67%  '

# Row 26 - SuiNetwork
Set-TextValue 'D26' '1.63'
$ws.Range('E26').Value = '  -5.36%  '

# Row 27 - Fetch.AI
$ws.Range('E27').Value = '  -1.95%  '

# Row 28 - Kaspa
$ws.Range('E28').Value = '  +1.74%  '

# Row 29 - Aptos
Set-TextValue 'D29' '8.17'
$ws.Range('E29').Value = '  -1.19%  '

# Row 30 - Binance-PegBSC-USD
$ws.Range('E30').Value = '  +0.19%  '

# Row 31 - Bittensor
Set-TextValue 'D31' '536.62'
$ws.Range('E31').Value = '  -1.14%  '

# Row 32 - PancakeSwap
Set-TextValue 'D32' '2.16'
$ws.Range('E32').Value = '  -2.94%  '

# Row 33 - ImmutableX
$ws.Range('E33').Value = '  -4.61%  '

# Row 34 - RenderToken
Set-TextValue 'D34' '6.55'
$ws.Range('E34').Value = '  +1.59%  '

# Row 35 - NEARProtocol
Set-TextValue 'D35' '5.40'
$ws.Range('E35').Value = '  -4.86%  '

# Row 36 - PolygonEcosystemToken
Set-TextValue 'D36' '0.426'
$ws.Range('E36').Value = '  -1.74%  '

# Row 37 - EthereumClassic
$ws.Range('E37').Value = '  -0.87%  '

# Row 38 - Monero
Set-TextValue 'D38' '160.53'
$ws.Range('E38').Value = '  -1.57%  '

# Row 39 - was Stacks, now FirstDigitalUSD (rows 39/40 swapped)
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D39' '1.00'
$ws.Range('E39').Value = '  +0.03%  '

# Row 40 - was FirstDigitalUSD, now Stacks (rows 39/40 swapped)
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D40' '1.97'
$ws.Range('E40').Value = '  -2.51%  '

# Row 41 - USDe
Set-TextValue 'D41' '0.999'

# Row 42 - OKB
Set-TextValue 'D42' '42.24'
$ws.Range('E42').Value = '  -0.88%  '

# Row 43 - Aave
Set-TextValue 'D43' '166.07'
$ws.Range('E43').Value = '  -0.36%  '

# Row 44 - Filecoin
Set-TextValue 'D44' '4.10'
$ws.Range('E44').Value = '  -2.36%  '

# Row 45 - Hedera
Set-TextValue 'D45' '0.0620'
$ws.Range('E45').Value = '  +0.17%  '

# Row 46 - InjectiveProtocol
Set-TextValue 'D46' '23.10'
$ws.Range('E46').Value = '  -0.50%  '

# Row 47 - dogwifhat
Set-TextValue 'D47' '2.23'
$ws.Range('E47').Value = '  -4.45%  '

# Row 48 - VeChain
Set-TextValue 'D48' '0.0263'
$ws.Range('E48').Value = '  -0.70%  '

# Row 49 - Mantle
Set-TextValue 'D49' '0.650'
$ws.Range('E49').Value = '  -1.42%  '

# Row 50 - EnergySwap
Set-TextValue 'D50' '20.34'

# Row 51 - Stellar
Set-TextValue 'D51' '0.0989'
$ws.Range('E51').Value = '  +0.21%  '
